$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "G3"  = 0.8
    "F6"  = 0.8272727272727273
    "F14" = 0.5864661654135339
    "F15" = 0.5714285714285714
    "G15" = 0.6962406015037594
    "F16" = 0.5699248120300752
    "G16" = 0.693233082706767
    "F17" = 0.5699248120300752
    "G17" = 0.693233082706767
    "G19" = 0.8833333333333333
    "F20" = 0.8
    "G20" = 0.8777777777777778
    "F21" = 0.7933333333333333
    "F22" = 0.8913043478260869
    "G24" = 0.9652173913043478
    "F25" = 0.8695652173913043
    "F26" = 0.52
    "G27" = 0.78
    "G29" = 0.7804999999999928
    "F30" = 0.8363636363636364
    "F42" = 0.5848484848484848
    "G42" = 0.7060606060606061
    "F43" = 0.5666666666666667
    "G43" = 0.693939393939394
    "F44" = 0.5696969696969697
    "F45" = 0.5681818181818181
    "G45" = 0.6909090909090909
    "F48" = 0.8631578947368421
    "F49" = 0.8578947368421053
    "G50" = 0.9134328358208955
    "F52" = 0.7313432835820896
    "G52" = 0.8716417910447761
    "G53" = 0.8716417910447761
    "F54" = 0.7893333333333333
    "F55" = 0.7573333333333333
    "G55" = 0.8853333333333333
    "F56" = 0.7626666666666666
    "G56" = 0.8853333333333333
    "G57" = 0.8746666666666667
    "G64" = 0.9619841269841258
    "F65" = 0.8476190476190476
    "F70" = 0.7943844492440605
    "G70" = 0.849244060475162
    "F71" = 0.7693304535637149
    "G71" = 0.827645788336933
    "G72" = 0.8250539956803455
    "F73" = 0.7649999999999999
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
